$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 49628.60690140428
$ws.Range("D2").Value = 45757427301.13125
$ws.Range("G2").Value = 40250401022.62417

$ws.Range("B3").Value = 50154.90453972395
$ws.Range("D3").Value = 45749326080.48063
$ws.Range("G3").Value = 38722355234.40031

$ws.Range("B4").Value = 54883.35478249052
$ws.Range("D4").Value = 45751322659.84856
$ws.Range("G4").Value = 40211298199.30566

$ws.Range("B5").Value = 55477.87105312177
$ws.Range("D5").Value = 45742452234.75478
$ws.Range("G5").Value = 38605788615.71039

$ws.Range("B6").Value = 43729.98641555249
$ws.Range("D6").Value = 44563079493.53232
$ws.Range("G6").Value = -9866133179.838198

$ws.Range("B7").Value = 44149.44968066624
$ws.Range("D7").Value = 44555183856.1448
$ws.Range("G7").Value = -11293786123.58261

$ws.Range("B8").Value = 43776.79034003732
$ws.Range("D8").Value = 44691000141.40578
$ws.Range("G8").Value = -9807838472.397394

$ws.Range("B9").Value = 44199.50623185693
$ws.Range("D9").Value = 44683170051.29977
$ws.Range("G9").Value = -11250772566.76314
